$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the crypto price/volume table (GitHub Actions data refresh).
# Price cells (column D) hold text like "47.229.79" / "1.00" that looks
# numeric, so NumberFormat is forced to Text ("@") right before each
# assignment to stop Excel from silently re-typing them as numbers
# (which would also strip formatting, e.g. "9.00" -> 9).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "47.229.79"
$ws.Range("E2").Value = "  +4.26%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.494.01"
$ws.Range("E3").Value = "  +2.84%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "323.22"
$ws.Range("E5").Value = "  +0.88%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "108.11"
$ws.Range("E6").Value = "  +4.57%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.525"
$ws.Range("E7").Value = "  +1.92%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.544"
$ws.Range("E9").Value = "  +2.84%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.96"
$ws.Range("E10").Value = "  +9.31%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0815"
$ws.Range("E11").Value = "  +1.84%  "
$ws.Range("E12").Value = "  +1.35%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.45"
$ws.Range("E13").Value = "  +1.30%  "
$ws.Range("E14").Value = "  +3.39%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.881.74"
$ws.Range("E15").Value = "  +2.76%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.491.82"
$ws.Range("E16").Value = "  +2.13%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.860"
$ws.Range("E17").Value = "  +3.31%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "47.172.14"
$ws.Range("E18").Value = "  +4.28%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.79"
$ws.Range("E19").Value = "  +4.30%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.67"
$ws.Range("E20").Value = "  +4.84%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0945"
$ws.Range("E21").Value = "  +2.05%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "70.95"
$ws.Range("E22").Value = "  -0.13%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.51"
$ws.Range("E23").Value = "  +5.00%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "250.28"
$ws.Range("E24").Value = "  +1.85%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.59"
$ws.Range("E25").Value = "  +3.78%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.23"
$ws.Range("E26").Value = "  +2.23%  "
$ws.Range("E27").Value = "  -0.03%  "
$ws.Range("B28").Value = "Cosmos"
$ws.Range("C28").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.02"
$ws.Range("E28").Value = "  +4.08%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.21"
$ws.Range("E29").Value = "  -3.20%  "
$ws.Range("B30").Value = "InjectiveProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "35.29"
$ws.Range("E30").Value = "  +5.55%  "
$ws.Range("B31").Value = "Kaspa"
$ws.Range("C31").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.136"
$ws.Range("E31").Value = "  +6.99%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "49.47"
$ws.Range("E32").Value = "  +0.74%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.96"
$ws.Range("E33").Value = "  -1.80%  "
$ws.Range("E34").Value = "  +3.68%  "
$ws.Range("E35").Value = "  +4.71%  "
$ws.Range("E36").Value = "  +0.18%  "
$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.76"
$ws.Range("E37").Value = "  +5.66%  "
$ws.Range("B38").Value = "ARBITRUM"
$ws.Range("C38").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.99"
$ws.Range("E38").Value = "  +5.81%  "
$ws.Range("E39").Value = "  +3.17%  "
$ws.Range("E40").Value = "  +1.62%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.24"
$ws.Range("E41").Value = "  -1.61%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "121.31"
$ws.Range("E42").Value = "  -5.66%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "21.53"
$ws.Range("E43").Value = "  +3.77%  "
$ws.Range("E44").Value = "  +2.70%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.975.51"
$ws.Range("E45").Value = "  +1.37%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.01"
$ws.Range("E46").Value = "  +2.30%  "
$ws.Range("E47").Value = "  -1.83%  "
$ws.Range("E48").Value = "  -0.13%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.02"
$ws.Range("E49").Value = "  -1.23%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.31"
$ws.Range("E50").Value = "  +10.67%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "78.43"
$ws.Range("E51").Value = "  +2.61%  "
